$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at E (shifts existing E:I data to F:J)
$ws.Columns("E").Insert()

# Give the new column a header / variable name
$ws.Range("E1").Value = "PressureTransducerSiteName"

# Populate PressureTransducerSiteName for the relevant stationary-antenna rows
$ws.Range("E7").Value = "Red Barn"
$ws.Range("E8").Value = "Red Barn"
$ws.Range("E9").Value = "Hitching Post"
$ws.Range("E10").Value = "Hitching Post"
$ws.Range("E11").Value = "Confluence"
$ws.Range("E12").Value = "Confluence"

# Give the new column a defined width (close to column D's width)
$ws.Columns("E").ColumnWidth = 15

# Update the active selection to reflect where the edit ended up
$ws.Range("C12").Select()
